$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 2987
$wsExhibit.Range("F12").Value = 7567
$wsExhibit.Range("F19").Value = 9253
$wsExhibit.Range("F36").Value = 1840
$wsExhibit.Range("F38").Value = 775
$wsExhibit.Range("F45").Value = 249

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value = 143

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2987
$wsAll.Range("F8").Value = 143
$wsAll.Range("F18").Value = 7567
$wsAll.Range("F24").Value = 9253
$wsAll.Range("F36").Value = 1842
$wsAll.Range("F38").Value = 775
$wsAll.Range("F46").Value = 249
